$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.658.08'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '2.475.65'
$ws.Range('E3').Value = '  +0.70%  '
$c = $ws.Range('D4')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E4').Value = '  -0.08%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '319.26'
$c.Style = 'Normal'
$ws.Range('E5').Value = '  +1.45%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '92.54'
$c.Style = 'Normal'
$ws.Range('E6').Value = '  +1.41%  '
$ws.Range('E7').Value = '  +0.79%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  +1.25%  '
$ws.Range('E10').Value = '  +8.96%  '
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '33.24'
$c.Style = 'Normal'
$ws.Range('E11').Value = '  +2.53%  '
$ws.Range('E12').Value = '  -0.18%  '
$ws.Range('D13').Value = '2.855.34'
$ws.Range('E13').Value = '  +0.58%  '
$ws.Range('E14').Value = '  +1.01%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '15.55'
$c.Style = 'Normal'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').Value = '2.472.45'
$ws.Range('E16').Value = '  +0.39%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.792'
$c.Style = 'Normal'
$ws.Range('E17').Value = '  +3.01%  '
$ws.Range('D18').Value = '41.602.54'
$ws.Range('E18').Value = '  +0.13%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '6.46'
$c.Style = 'Normal'
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = '0.0₃0943'
$ws.Range('E20').Value = '  +0.86%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '70.70'
$c.Style = 'Normal'
$ws.Range('E21').Value = '  -0.32%  '
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '11.30'
$c.Style = 'Normal'
$ws.Range('E22').Value = '  +0.15%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '240.41'
$c.Style = 'Normal'
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('E24').Value = '  +1.55%  '
$ws.Range('E25').Value = '  +2.64%  '
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.Style = 'Normal'
$ws.Range('E26').Value = '  +0.04%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '24.86'
$c.Style = 'Normal'
$ws.Range('E27').Value = '  +2.52%  '
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '2.24'
$c.Style = 'Normal'
$ws.Range('E28').Value = '  -0.44%  '
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '9.71'
$c.Style = 'Normal'
$ws.Range('E29').Value = '  +0.70%  '
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '36.68'
$c.Style = 'Normal'
$ws.Range('E30').Value = '  +4.86%  '
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '157.09'
$c.Style = 'Normal'
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('E32').Value = '  +0.94%  '
$ws.Range('E33').Value = '  -0.15%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.0766'
$c.Style = 'Normal'
$ws.Range('E34').Value = '  +1.17%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.55'
$c.Style = 'Normal'
$ws.Range('E35').Value = '  -0.86%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '17.27'
$c.Style = 'Normal'
$ws.Range('E36').Value = '  -0.45%  '
$ws.Range('E37').Value = '  +4.86%  '
$ws.Range('E38').Value = '  +1.65%  '
$ws.Range('B39').Value = 'Kaspa'
$ws.Range('C39').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.105'
$c.Style = 'Normal'
$ws.Range('E39').Value = '  +2.77%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.90'
$c.Style = 'Normal'
$ws.Range('E40').Value = '  +1.03%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '4.01'
$c.Style = 'Normal'
$ws.Range('E41').Value = '  +1.89%  '
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '2.49'
$c.Style = 'Normal'
$ws.Range('E42').Value = '  +3.49%  '
$ws.Range('D43').Value = '1.989.37'
$ws.Range('E43').Value = '  +1.13%  '
$ws.Range('E44').Value = '  +0.52%  '
$ws.Range('E45').Value = '  +1.82%  '
$ws.Range('E46').Value = '  +2.95%  '
$ws.Range('E47').Value = '  +5.88%  '
$ws.Range('D48').Value = '2.712.80'
$ws.Range('E48').Value = '  +0.60%  '
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '97.96'
$c.Style = 'Normal'
$ws.Range('E49').Value = '  +1.79%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '75.75'
$c.Style = 'Normal'
$ws.Range('E50').Value = '  +5.63%  '
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '67.18'
$c.Style = 'Normal'
$ws.Range('E51').Value = '  +1.48%  '
